$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D from Excel auto-converting numeric-looking text into numbers,
# so values like "246.40" or "5.41" remain stored as text (matching the source data).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.065.25"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "2.051.09"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "246.40"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "0.660"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").Value = "58.52"
$ws.Range("E7").Value = "  -4.87%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.379"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("E11").Value = "  +2.39%  "
$ws.Range("D12").Value = "15.42"
$ws.Range("E12").Value = "  -5.48%  "
$ws.Range("E13").Value = "  +7.80%  "
$ws.Range("D14").Value = "2.346.97"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D16").Value = "2.040.63"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "18.31"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "37.011.20"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "73.86"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "0.0₃0888"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").Value = "5.41"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").Value = "238.83"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("D26").Value = "168.73"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("D28").Value = "20.03"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "5.59"
$ws.Range("E29").Value = "  +15.99%  "
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "1.13"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  +4.31%  "
$ws.Range("D33").Value = "0.0615"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  +6.31%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "0.0843"
$ws.Range("E37").Value = "  -5.69%  "
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").Value = "5.25"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "0.0983"
$ws.Range("E41").Value = "  -9.71%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0223"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").Value = "97.89"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "17.02"
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.301.87"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").Value = "  -4.93%  "
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D49").Value = "6.76"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").Value = "3.64"
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").Value = "2.229.49"
$ws.Range("E51").Value = "  -0.52%  "

# Restore default (Normal) style on column D so no stray number-format styling remains
$ws.Range("D2:D51").Style = "Normal"
